$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose text values look like plain numbers need to be forced to Text
# format first, otherwise Excel will silently convert them into floating point
# numbers (losing the exact original textual representation).
$textCells = @("D5","D6","D9","D10","D11","D12","D13","D14","D16","D20","D21","D23","D27","D28","D31","D32","D33","D34","D36","D39","D40","D41","D43","D44","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '43.828.25'
$ws.Range("E2").Value = '  +0.44%  '
$ws.Range("D3").Value = '2.297.30'
$ws.Range("E3").Value = '  +0.26%  '
$ws.Range("E4").Value = '  +0.28%  '
$ws.Range("D5").Value = '114.35'
$ws.Range("E5").Value = '  +18.68%  '
$ws.Range("D6").Value = '268.75'
$ws.Range("E6").Value = '  +0.28%  '
$ws.Range("E7").Value = '  +1.69%  '
$ws.Range("E8").Value = '  +0.32%  '
$ws.Range("D9").Value = '0.623'
$ws.Range("E9").Value = '  +1.97%  '
$ws.Range("D10").Value = '48.42'
$ws.Range("E10").Value = '  +5.39%  '
$ws.Range("D11").Value = '0.0953'
$ws.Range("E11").Value = '  +1.79%  '
$ws.Range("D12").Value = '8.93'
$ws.Range("E12").Value = '  +13.13%  '
$ws.Range("D13").Value = '0.107'
$ws.Range("E13").Value = '  +1.43%  '
$ws.Range("D14").Value = '15.74'
$ws.Range("E14").Value = '  +3.71%  '
$ws.Range("D15").Value = '2.639.76'
$ws.Range("E15").Value = '  +0.29%  '
$ws.Range("D16").Value = '0.851'
$ws.Range("E16").Value = '  +0.67%  '
$ws.Range("D17").Value = '2.296.99'
$ws.Range("E17").Value = '  +0.22%  '
$ws.Range("D18").Value = '43.744.86'
$ws.Range("E18").Value = '  +0.36%  '
$ws.Range("E19").Value = '  +2.17%  '
$ws.Range("D20").Value = '6.63'
$ws.Range("E20").Value = '  +6.90%  '
$ws.Range("D21").Value = '72.58'
$ws.Range("E21").Value = '  +0.43%  '
$ws.Range("E22").Value = '  -0.24%  '
$ws.Range("D23").Value = '233.25'
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("E24").Value = '  +7.31%  '
$ws.Range("E25").Value = '  +11.74%  '
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("D27").Value = '11.69'
$ws.Range("E27").Value = '  +4.14%  '
$ws.Range("D28").Value = '41.93'
$ws.Range("E28").Value = '  +4.41%  '
$ws.Range("E29").Value = '  -2.04%  '
$ws.Range("E30").Value = '  +1.19%  '
$ws.Range("D31").Value = '176.92'
$ws.Range("E31").Value = '  +0.73%  '
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").Value = '0.0934'
$ws.Range("E32").Value = '  +4.41%  '
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").Value = '21.69'
$ws.Range("E33").Value = '  -0.81%  '
$ws.Range("D34").Value = '5.64'
$ws.Range("E34").Value = '  +5.15%  '
$ws.Range("E35").Value = '  +1.14%  '
$ws.Range("D36").Value = '4.76'
$ws.Range("E36").Value = '  +9.17%  '
$ws.Range("E37").Value = '  +2.62%  '
$ws.Range("E38").Value = '  +0.36%  '
$ws.Range("D39").Value = '3.82'
$ws.Range("E39").Value = '  +12.21%  '
$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D40").Value = '0.246'
$ws.Range("E40").Value = '  +0.30%  '
$ws.Range("B41").Value = 'Celestia'
$ws.Range("C41").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D41").Value = '13.98'
$ws.Range("E41").Value = '  +13.88%  '
$ws.Range("E42").Value = '  +4.19%  '
$ws.Range("D43").Value = '72.05'
$ws.Range("E43").Value = '  +11.50%  '
$ws.Range("D44").Value = '6.21'
$ws.Range("E44").Value = '  +19.54%  '
$ws.Range("E45").Value = '  +4.24%  '
$ws.Range("E46").Value = '  +0.23%  '
$ws.Range("E47").Value = '  +0.33%  '
$ws.Range("E48").Value = '  -0.82%  '
$ws.Range("D49").Value = '102.57'
$ws.Range("E49").Value = '  +5.27%  '
$ws.Range("D50").Value = '1.23'
$ws.Range("E50").Value = '  +3.66%  '
$ws.Range("D51").Value = '0.449'
$ws.Range("E51").Value = '  +6.37%  '
